$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A14").Value = 13
$ws.Range("B14").Value = "Friday, Jan 13"
$ws.Range("C14").Value = "9:05 AM"
$ws.Range("D14").Value = "W95177"
$ws.Range("E14").Value = "London"
$ws.Range("F14").Value = "(LTN)"
$ws.Range("G14").Value = "Wizz Air "
$ws.Range("H14").Value = "A320"
$ws.Range("I14").Value = "(G-WUKF)"
$ws.Range("J14").Value = "8:41 AM"
$ws.Range("K14").Borders.LineStyle = -4142
$ws.Range("L14").Value = "0 hours, -24 minutes"
$ws.Range("M14").Borders.LineStyle = -4142
